# The deck currently ships two theme parts:
#   ppt/theme/theme1.xml -> "Integral"    (used by the Slide Master)
#   ppt/theme/theme2.xml -> "Office Theme" (used by the Notes Master)
# The commit swaps them: the Slide Master's design becomes the plain
# "Office Theme" colour palette, while the old "Integral" palette is
# kept around as the secondary theme. Re-colour the presentation's
# (Slide Master's) theme color scheme to the stock "Office Theme"
# values via the ThemeColorScheme object, which is the object-model
# surface PowerPoint exposes for editing a design's 12 theme colors.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# MsoThemeColorSchemeIndex order: 1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
# 5-10 Accent1-6, 11 Hyperlink, 12 FollowedHyperlink.
$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
